# Mark additional "Users" and "Network" methods as Done on the Methods sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Methods")

# Finished user methods (rows 71-77)
$doneRows = 71,72,73,74,75,76,77,81,82,83,84

foreach ($r in $doneRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "Done"
    $cell.Style = "Good"
}

# Move the active selection to A85 (last edit position)
$ws.Activate()
$ws.Range("A85").Select()
